$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.596.00'
$ws.Range('E2').Value = '  +2.94%  '
$ws.Range('D3').Value = '2.353.02'
$ws.Range('E3').Value = '  +5.95%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''314.14'
$ws.Range('E5').Value = '  +5.81%  '
$ws.Range('D6').Value = '''110.31'
$ws.Range('E6').Value = '  +1.99%  '
$ws.Range('D7').Value = '''0.644'
$ws.Range('E7').Value = '  +2.98%  '
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').Value = '''0.636'
$ws.Range('E9').Value = '  +6.29%  '
$ws.Range('D10').Value = '''43.45'
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('E11').Value = '  +3.00%  '
$ws.Range('D12').Value = '''8.87'
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').Value = '''1.04'
$ws.Range('E13').Value = '  +4.49%  '
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D15').Value = '''16.39'
$ws.Range('E15').Value = '  +9.37%  '
$ws.Range('D16').Value = '2.708.53'
$ws.Range('E16').Value = '  +6.18%  '
$ws.Range('D17').Value = '2.423.55'
$ws.Range('E17').Value = '  +8.50%  '
$ws.Range('D18').Value = '43.582.92'
$ws.Range('E18').Value = '  +3.07%  '
$ws.Range('E19').Value = '  +3.56%  '
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('D21').Value = '''75.78'
$ws.Range('E21').Value = '  +4.60%  '
$ws.Range('D22').Value = '''3.46'
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').Value = '''2.58'
$ws.Range('E23').Value = '  +11.87%  '
$ws.Range('D24').Value = '''256.83'
$ws.Range('E24').Value = '  +12.39%  '
$ws.Range('D25').Value = '''9.15'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').Value = '''12.09'
$ws.Range('E26').Value = '  +4.40%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '''39.32'
$ws.Range('E28').Value = '  +3.01%  '
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('D30').Value = '''22.44'
$ws.Range('E30').Value = '  +7.15%  '
$ws.Range('D31').Value = '''174.13'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('E33').Value = '  +3.87%  '
$ws.Range('E34').Value = '  +7.83%  '
$ws.Range('E35').Value = '  +5.49%  '
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('D37').Value = '''4.18'
$ws.Range('E37').Value = '  -4.22%  '
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('E39').Value = '  +2.35%  '
$ws.Range('D40').Value = '''2.72'
$ws.Range('E40').Value = '  +12.51%  '
$ws.Range('D41').Value = '''72.77'
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('E42').Value = '  +14.24%  '
$ws.Range('E43').Value = '  +1.42%  '
$ws.Range('D44').Value = '''12.87'
$ws.Range('E44').Value = '  +1.99%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = '''5.67'
$ws.Range('E46').Value = '  +4.84%  '
$ws.Range('D47').Value = '''9.32'
$ws.Range('E47').Value = '  +10.95%  '
$ws.Range('D48').Value = '''111.63'
$ws.Range('E48').Value = '  +8.01%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('E50').Value = '  +2.96%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '''0.462'
$ws.Range('E51').Value = '  +5.31%  '
